# KPI.xlsx update: Target headcount moved from 720 to 810 for every row,
# and several "percent complete" drivers (G column) were retuned to match
# the new target. K/L columns are formula-driven off G/H and recompute
# automatically; only the selected cell moves (K11 -> G29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- H2:H37 -> 810 (was 720 for every row) ---
for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 8).Value = 810
}

# --- G3 is a literal (not a formula) that tracked H3 -> update together ---
$ws.Range("G3").Value = 810

# --- G column formulas: percentage-of-target drivers that changed ---
$ws.Range("G6").Formula  = "=49%*H6"
$ws.Range("G7").Formula  = "=82%*H7"
$ws.Range("G13").Formula = "=60%*H13"
$ws.Range("G14").Formula = "=16%*H14"
$ws.Range("G16").Formula = "=30%*H16"
$ws.Range("G17").Formula = "=99%*H17"
$ws.Range("G18").Formula = "=88%*H18"
$ws.Range("G21").Formula = "=43%*H21"
$ws.Range("G22").Formula = "=67%*H22"
$ws.Range("G27").Formula = "=23%*H27"
$ws.Range("G29").Formula = "=34%*H29"
$ws.Range("G31").Formula = "=69%*H31"
$ws.Range("G32").Formula = "=40%*H32"
$ws.Range("G35").Formula = "=32%*H35"
$ws.Range("G36").Formula = "=66%*H36"

# --- Selection moved from K11 to G29 ---
[void]$ws.Range("G29").Select()

Write-Output "done"
